# Add three new webcam/location rows (40, 41, 42) to the "location-1" sheet,
# matching the data that was appended to the source workbook.
#
# Calculation is switched to manual first so that re-touching the worksheet
# does not force a recalculation of the existing IsYouTubeVideoValid(...)
# shared-formula column (that UDF is not available in this workbook, and a
# recalculation would blow away the cached TRUE/FALSE results already
# stored for rows 2-39).
$excel.Calculation = -4135

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 40 - Agdao Farmers Market (Davao, Philippines)
# ---------------------------------------------------------------------
$ws.Range("A40").Value = "LIVE, SHOPPING, MARKET"
$ws.Range("B40").Value = "7.082030630801781, 125.6236221942979"
$ws.Range("C40").Value = "Agdao Farmers Market, Davao City Cam 1"
$ws.Range("D40").Value = "Davao"
$ws.Range("E40").Value = "Philippines"
$ws.Range("F40").Value = "mPqGLVpYN5Q"
$ws.Range("G40").Formula = "=IsYouTubeVideoValid(F40)"

# Formatting: same pattern as the "plain" data rows (e.g. row 27):
#   A/C/D/E boxed (thin left/right border), B/F unformatted, G boxed.
$ws.Range("A27:G27").Copy()
$ws.Range("A40:G40").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Row 41 - PTZ Live Camera in Philippines, Construction & Market
# ---------------------------------------------------------------------
$ws.Range("A41").Value = "LIVE, TRAFFIC"
$ws.Range("B41").Value = "7.081855910707648, 125.62407457207543"
$ws.Range("C41").Value = "PTZ Live Camera in Philippines, Construction & Market"
$ws.Range("D41").Value = "Davao"
$ws.Range("E41").Value = "Philippines"
$ws.Range("F41").Value = "t45_gP7I82I"
$ws.Range("G41").Formula = "=IsYouTubeVideoValid(F41)"

# Same formatting pattern as row 40.
$ws.Range("A27:G27").Copy()
$ws.Range("A41:G41").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Row 42 - Live Hummingbird Cam in Peru
# ---------------------------------------------------------------------
$ws.Range("A42").Value = "LIVE, BIRD, NATURE"
$ws.Range("B42").Value = "-13.256655478925458, -72.17352170709896"
$ws.Range("C42").Value = "Live Hummingbird Cam in Peru"
$ws.Range("D42").Value = "av.pucara"
$ws.Range("E42").Value = "Peru"
$ws.Range("F42").Value = "Ej_bvcIlByY"

# Formatting follows the pattern used for rows 38/39: A/C/D/E/F boxed,
# B uses the quote-prefixed style (its text begins with "-", like a
# negative number). There is no YouTube-validity formula for this row.
$ws.Range("A38:E38").Copy()
$ws.Range("A42:E42").PasteSpecial(-4122)
$ws.Range("C38").Copy()
$ws.Range("F42").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Restore automatic calculation and update the view/selection to match
# the state the workbook was left in.
# ---------------------------------------------------------------------
$excel.Calculation = -4105
$ws.Range("F45").Select()
